# Se actualizan valores de inventario
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Inventario principal" ---
$ws1 = $wb.Worksheets.Item("Inventario principal")

# Update inventory quantities
$ws1.Range("D2").Value = 10
$ws1.Range("D4").Value = 15

# Update the active selection to G11
$ws1.Range("G11").Select()

# --- Sheet 2: "Reporte 2025-04-14" ---
$ws2 = $wb.Worksheets.Item("Reporte 2025-04-14")

# Update inventory quantities (mirrors sheet 1)
$ws2.Range("D2").Value = 10
$ws2.Range("D4").Value = 15
